$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "23.191.73"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.37%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.602.27"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.02%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "1.000"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.03%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "303.35"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.63%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3783"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.16%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "52.05"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +4.37%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3619"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.87%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.272"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.42%  "
$ws.Range("E11").Value = "  -0.03%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08120"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.41%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "22.78"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.77%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.604"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.18%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.421"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.43%  "
$ws.Range("E16").Value = "  -1.28%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.603.35"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.24%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "94.03"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.27%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06878"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "18.08"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.94%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.550"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.68%  "
$ws.Range("E22").Value = "  -0.01%  "
$ws.Range("E23").Value = "  -0.53%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "23.189.89"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.32%  "
$ws.Range("E25").Value = "  +2.02%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.975"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +8.21%  "
$ws.Range("E27").Value = "  +0.54%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "149.42"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.23%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.251"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.25%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "133.95"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.96%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.363"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.39%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.780"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.11%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.779.33"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.01%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.9714"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.67%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.07521"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.34%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "10.30"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.00%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02721"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.13%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2505"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.89%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.08803"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.11%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.078"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -3.14%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.7114"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.31%  "
$ws.Range("E42").Value = "  -0.73%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "12.49"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.93%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "15.66"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.05%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.6538"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.13%  "
$ws.Range("E46").Value = "  -0.21%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.015"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.49%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "132.39"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.19%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.07955"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.18%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.202"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.10%  "
$ws.Range("E51").Value = "  +1.14%  "
